$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.775.91'
$ws.Range("E2").Value = '  -1.47%  '

# Row 3
$ws.Range("D3").Value = '1.548.94'

# Row 4
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").Value = '''206.22'
$ws.Range("E5").Value = '  -0.60%  '

# Row 6
$ws.Range("E6").Value = '  -1.90%  '

# Row 7
$ws.Range("E7").Value = '  -0.16%  '

# Row 8
$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").Value = '''21.44'
$ws.Range("E8").Value = '  -3.86%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.245'
$ws.Range("E9").Value = '  -1.41%  '

# Row 10
$ws.Range("E10").Value = '  -1.34%  '

# Row 11
$ws.Range("E11").Value = '  -1.57%  '

# Row 12
$ws.Range("D12").Value = '1.769.11'
$ws.Range("E12").Value = '  -1.71%  '

# Row 13
$ws.Range("D13").Value = '1.559.87'
$ws.Range("E13").Value = '  -1.09%  '

# Row 14
$ws.Range("E14").Value = '  -2.70%  '

# Row 15
$ws.Range("D15").Value = '''0.512'
$ws.Range("E15").Value = '  -1.35%  '

# Row 16
$ws.Range("D16").Value = '26.773.00'

# Row 17
$ws.Range("D17").Value = '''61.15'
$ws.Range("E17").Value = '  -1.89%  '

# Row 18
$ws.Range("D18").Value = '''213.19'
$ws.Range("E18").Value = '  -0.66%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0687'
$ws.Range("E19").Value = '  +0.20%  '

# Row 20
$ws.Range("E20").Value = '  -2.07%  '

# Row 21
$ws.Range("E21").Value = '  -0.02%  '

# Row 22
$ws.Range("E22").Value = '  -1.07%  '

# Row 23
$ws.Range("D23").Value = '''8.95'
$ws.Range("E23").Value = '  -5.62%  '

# Row 24
$ws.Range("D24").Value = '''1.99'
$ws.Range("E24").Value = '  -1.56%  '

# Row 25
$ws.Range("D25").Value = '''153.12'
$ws.Range("E25").Value = '  +0.38%  '

# Row 26
$ws.Range("D26").Value = '''6.52'
$ws.Range("E26").Value = '  -2.74%  '

# Row 27
$ws.Range("D27").Value = '''14.88'
$ws.Range("E27").Value = '  -0.50%  '

# Row 28
$ws.Range("E28").Value = '  -0.07%  '

# Row 29
$ws.Range("E29").Value = '  -1.52%  '

# Row 30
$ws.Range("E30").Value = '  -0.86%  '

# Row 31
$ws.Range("E31").Value = '  -1.35%  '

# Row 32
$ws.Range("D32").Value = '''3.19'
$ws.Range("E32").Value = '  +0.23%  '

# Row 33
$ws.Range("D33").Value = '1.344.86'
$ws.Range("E33").Value = '  -3.89%  '

# Row 34
$ws.Range("E34").Value = '  -0.26%  '

# Row 35
$ws.Range("E35").Value = '  -3.55%  '

# Row 36
$ws.Range("D36").Value = '''2.28'
$ws.Range("E36").Value = '  -0.48%  '

# Row 37
$ws.Range("D37").Value = '''0.928'
$ws.Range("E37").Value = '  -1.16%  '

# Row 38
$ws.Range("D38").Value = '''0.0164'
$ws.Range("E38").Value = '  -0.80%  '

# Row 39
$ws.Range("D39").Value = '''0.521'
$ws.Range("E39").Value = '  +0.86%  '

# Row 40
$ws.Range("D40").Value = '''0.802'
$ws.Range("E40").Value = '  -1.85%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''5.71'
$ws.Range("E41").Value = '  +4.89%  '

# Row 42
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '''0.992'
$ws.Range("E42").Value = '  -1.28%  '

# Row 43
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.19'
$ws.Range("E43").Value = '  -0.02%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''1.76'
$ws.Range("E44").Value = '  -4.71%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''62.84'
$ws.Range("E45").Value = '  -1.58%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.682.76'
$ws.Range("E46").Value = '  -1.68%  '

# Row 47
$ws.Range("B47").Value = 'mCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D47").Value = '''2.24'
$ws.Range("E47").Value = '  -4.29%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''85.89'
$ws.Range("E48").Value = '  +0.12%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.0516'
$ws.Range("E49").Value = '  +4.36%  '

# Row 50
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₇0975'
$ws.Range("E50").Value = '  -0.28%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.0950'
$ws.Range("E51").Value = '  -0.19%  '
